$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(203957296, Omri Ben Shabat: 5,-2)"
$ws.Range("B1").Value = "(206532695, Matan Vakrat: -7,7)"
$ws.Range("C1").Value = "(302962915, Asher  Odeh: -1,9)"
$ws.Range("D1").Value = "(308035542, Anastasia  Kubi: -6,2)"
$ws.Range("E1").Value = "(311177802, Christina  Uksusman: 0,9)"
$ws.Range("F1").Value = "(305251175, Or  Leder: -4,1)"

$ws.Range("A3").Value = "cost: 409.72557563737143"
$ws.Range("A4").Value = "time: 48.09069695467143"
